$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B:E are treated as text so numeric-looking strings
# (e.g. "260.10", "0.111", prices with thousand-dot separators, percentages)
# are preserved exactly as authored instead of being coerced into numbers.
$ws.Range("B2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '43.832.21'
$ws.Range("E2").Value = '  +2.18%  '

# Row 3
$ws.Range("D3").Value = '2.215.96'
$ws.Range("E3").Value = '  +0.36%  '

# Row 4
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$ws.Range("D5").Value = '260.10'
$ws.Range("E5").Value = '  +2.25%  '

# Row 6
$ws.Range("D6").Value = '86.68'
$ws.Range("E6").Value = '  +13.89%  '

# Row 7
$ws.Range("D7").Value = '0.619'
$ws.Range("E7").Value = '  +0.64%  '

# Row 8
$ws.Range("E8").Value = '  +0.05%  '

# Row 9
$ws.Range("D9").Value = '0.598'
$ws.Range("E9").Value = '  +0.95%  '

# Row 10
$ws.Range("D10").Value = '45.52'
$ws.Range("E10").Value = '  +8.40%  '

# Row 11
$ws.Range("D11").Value = '0.0920'
$ws.Range("E11").Value = '  +1.21%  '

# Row 12
$ws.Range("E12").Value = '  +7.94%  '

# Row 13
$ws.Range("E13").Value = '  +1.75%  '

# Row 14
$ws.Range("D14").Value = '2.545.18'
$ws.Range("E14").Value = '  +0.22%  '

# Row 15
$ws.Range("D15").Value = '14.47'
$ws.Range("E15").Value = '  +0.43%  '

# Row 16
$ws.Range("D16").Value = '2.207.41'
$ws.Range("E16").Value = '  +0.01%  '

# Row 17
$ws.Range("E17").Value = '  +0.61%  '

# Row 18
$ws.Range("D18").Value = '43.780.01'
$ws.Range("E18").Value = '  +2.32%  '

# Row 19
$ws.Range("D19").Value = '0.0000104'
$ws.Range("E19").Value = '  +1.01%  '

# Row 20
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '5.94'
$ws.Range("E20").Value = '  +0.31%  '

# Row 21
$ws.Range("B21").Value = 'Litecoin'
$ws.Range("C21").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D21").Value = '69.85'
$ws.Range("E21").Value = '  -1.89%  '

# Row 22
$ws.Range("E22").Value = '  +7.67%  '

# Row 23
$ws.Range("D23").Value = '231.39'
$ws.Range("E23").Value = '  +0.97%  '

# Row 24
$ws.Range("D24").Value = '8.97'
$ws.Range("E24").Value = '  -3.05%  '

# Row 25
$ws.Range("E25").Value = '  +0.04%  '

# Row 26
$ws.Range("D26").Value = '3.54'
$ws.Range("E26").Value = '  +5.47%  '

# Row 27
$ws.Range("E27").Value = '  +0.44%  '

# Row 28
$ws.Range("D28").Value = '39.64'
$ws.Range("E28").Value = '  +1.00%  '

# Row 29
$ws.Range("D29").Value = '2.25'
$ws.Range("E29").Value = '  +2.88%  '

# Row 30
$ws.Range("E30").Value = '  +2.23%  '

# Row 31
$ws.Range("D31").Value = '174.12'
$ws.Range("E31").Value = '  +0.49%  '

# Row 32
$ws.Range("D32").Value = '20.47'
$ws.Range("E32").Value = '  +1.39%  '

# Row 33
$ws.Range("E33").Value = '  +1.48%  '

# Row 34
$ws.Range("D34").Value = '5.37'
$ws.Range("E34").Value = '  +3.44%  '

# Row 35
$ws.Range("E35").Value = '  +1.82%  '

# Row 36
$ws.Range("D36").Value = '0.111'
$ws.Range("E36").Value = '  +3.16%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '4.49'
$ws.Range("E37").Value = '  +4.69%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.0356'
$ws.Range("E38").Value = '  +0.86%  '

# Row 39
$ws.Range("D39").Value = '12.65'
$ws.Range("E39").Value = '  +1.69%  '

# Row 40
$ws.Range("D40").Value = '2.90'
$ws.Range("E40").Value = '  +6.03%  '

# Row 41
$ws.Range("E41").Value = '  +0.06%  '

# Row 42
$ws.Range("D42").Value = '63.31'
$ws.Range("E42").Value = '  +5.81%  '

# Row 43
$ws.Range("D43").Value = '5.49'
$ws.Range("E43").Value = '  +4.54%  '

# Row 44
$ws.Range("E44").Value = '  +0.70%  '

# Row 45
$ws.Range("D45").Value = '100.63'
$ws.Range("E45").Value = '  -0.82%  '

# Row 46
$ws.Range("D46").Value = '8.34'
$ws.Range("E46").Value = '  +0.74%  '

# Row 47
$ws.Range("E47").Value = '  +0.36%  '

# Row 48
$ws.Range("E48").Value = '  +4.53%  '

# Row 49
$ws.Range("E49").Value = '  +1.01%  '

# Row 50
$ws.Range("D50").Value = '0.435'
$ws.Range("E50").Value = '  -3.49%  '

# Row 51
$ws.Range("E51").Value = '  +4.84%  '
